$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18
$ws.Range("B3").Value = 31
$ws.Range("B4").Value = 16
$ws.Range("B5").Value = 22
$ws.Range("B6").Value = 21
